$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D3").Select()
